# redid some incorrect times on vectorization
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected vectorization timings (rows 15-20, column E) ---
$ws.Range("E15").Value = 3.043
$ws.Range("E16").Value = 2.436
$ws.Range("E17").Value = 1.665
$ws.Range("E18").Value = 1.528
$ws.Range("E19").Value = 3.163
$ws.Range("E20").Value = 1.558

# --- Row 30 / D30 gains the same wrapped-text style used by rows 26-29,
#     and the row grows tall enough to show the wrapped comment. ---
$ws.Range("D30").WrapText = $true
$ws.Rows(30).RowHeight = 34

# --- Update the visible window position / selection to match where the
#     author was last working. ---
$ws.Range("J19").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
